$wb = $excel.ActiveWorkbook

# Update Metadata sheet timestamp
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = "05 Nov 2025, 02:38 PM"

# Update "1 Year" column (F) values on Industry Analysis sheet
$iaWs = $wb.Worksheets.Item("Industry Analysis")
$iaWs.Range("F2").Value = 21.0016
$iaWs.Range("F3").Value = -16.2396
$iaWs.Range("F4").Value = 27.1317
$iaWs.Range("F5").Value = -50.6494
$iaWs.Range("F6").Value = 53.2813
$iaWs.Range("F7").Value = -8.106199999999999
$iaWs.Range("F8").Value = -9.552099999999999
$iaWs.Range("F9").Value = 36.3756
$iaWs.Range("F10").Value = -6.1314
$iaWs.Range("F11").Value = 31.9081
$iaWs.Range("F12").Value = -18.4955
$iaWs.Range("F13").Value = 14.0155
$iaWs.Range("F14").Value = -36.0718
$iaWs.Range("F15").Value = -0.1622
$iaWs.Range("F16").Value = 0.1459
$iaWs.Range("F17").Value = -22.0012
$iaWs.Range("F18").Value = 1.0561
$iaWs.Range("F19").Value = -27.708
$iaWs.Range("F20").Value = 47.7309
$iaWs.Range("F21").Value = 12.0959
$iaWs.Range("F22").Value = 95.1491
$iaWs.Range("F23").Value = -50.2657
$iaWs.Range("F24").Value = -13.3427
$iaWs.Range("F25").Value = -9.9316
$iaWs.Range("F26").Value = 5.8244
$iaWs.Range("F27").Value = -32.7692
$iaWs.Range("F28").Value = -24.8224
$iaWs.Range("F29").Value = -18.4191
$iaWs.Range("F30").Value = 25.8569
$iaWs.Range("F31").Value = 58.4712
$iaWs.Range("F32").Value = -3.3862
$iaWs.Range("F33").Value = -6.3282
$iaWs.Range("F34").Value = 27.7203
$iaWs.Range("F35").Value = 4.4873
$iaWs.Range("F36").Value = -4.9458
$iaWs.Range("F37").Value = 3.6074
$iaWs.Range("F38").Value = -23.3973
$iaWs.Range("F39").Value = 8.7355
$iaWs.Range("F40").Value = -5.8541
$iaWs.Range("F41").Value = -8.3934
$iaWs.Range("F42").Value = 20.3818
$iaWs.Range("F43").Value = 14.3164
$iaWs.Range("F44").Value = -12.6846
$iaWs.Range("F45").Value = 28.4075
$iaWs.Range("F46").Value = -1.1135
$iaWs.Range("F47").Value = -37.1997
$iaWs.Range("F48").Value = -29.8569
$iaWs.Range("F49").Value = -27.5511
$iaWs.Range("F50").Value = -49.7478
$iaWs.Range("F51").Value = -51.8002
$iaWs.Range("F52").Value = -38.5254
$iaWs.Range("F53").Value = -12.4886
$iaWs.Range("F54").Value = -5.0725
$iaWs.Range("F55").Value = -17.7445
$iaWs.Range("F56").Value = -26.636
$iaWs.Range("F57").Value = -29.3361
$iaWs.Range("F58").Value = -11.9574
$iaWs.Range("F59").Value = -24.5687
$iaWs.Range("F60").Value = -12.3
$iaWs.Range("F61").Value = -10.9446
$iaWs.Range("F62").Value = -17.1229
$iaWs.Range("F63").Value = -9.5038
$iaWs.Range("F64").Value = 54.2749
$iaWs.Range("F65").Value = -43.4736
$iaWs.Range("F66").Value = 13.2687
$iaWs.Range("F67").Value = 12.7149
$iaWs.Range("F68").Value = 24.8057
$iaWs.Range("F69").Value = -17.0328
$iaWs.Range("F70").Value = -6.8927
$iaWs.Range("F71").Value = 13.6034
$iaWs.Range("F72").Value = 3.9995
$iaWs.Range("F73").Value = -16.226
$iaWs.Range("F74").Value = -16.2448
$iaWs.Range("F75").Value = 28.6924
$iaWs.Range("F76").Value = 48.9752
